$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the format from H1
# so they match the existing bold/centered/bordered header style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) for rows 2-82: row, I-value, J-value
$data = @(
    @(2, 2, 2),
    @(3, 6, 6),
    @(4, 5, 6),
    @(5, 6, 6),
    @(6, 2, 2),
    @(7, 7, 8),
    @(8, 10, 10),
    @(9, 7, 7),
    @(10, 5, 6),
    @(11, 8, 9),
    @(12, 7, 7),
    @(13, 6, 6),
    @(14, 4, 5),
    @(15, 7, 7),
    @(16, 5, 6),
    @(17, 4, 5),
    @(18, 6, 6),
    @(19, 8, 8),
    @(20, 6, 7),
    @(21, 4, 4),
    @(22, 5, 6),
    @(23, 4, 5),
    @(24, 8, 8),
    @(25, 6, 6),
    @(26, 5, 5),
    @(27, 6, 6),
    @(28, 6, 7),
    @(29, 9, 9),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 4, 5),
    @(33, 7, 7),
    @(34, 6, 6),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 7, 8),
    @(38, 7, 7),
    @(39, 9, 9),
    @(40, 6, 6),
    @(41, 8, 8),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 8, 9),
    @(45, 6, 6),
    @(46, 6, 7),
    @(47, 7, 7),
    @(48, 7, 7),
    @(49, 9, 9),
    @(50, 8, 8),
    @(51, 5, 6),
    @(52, 5, 6),
    @(53, 6, 7),
    @(54, 8, 8),
    @(55, 5, 6),
    @(56, 8, 8),
    @(57, 7, 7),
    @(58, 8, 9),
    @(59, 5, 6),
    @(60, 7, 7),
    @(61, 7, 7),
    @(62, 5, 6),
    @(63, 5, 6),
    @(64, 8, 8),
    @(65, 4, 5),
    @(66, 7, 7),
    @(67, 8, 8),
    @(68, 7, 7),
    @(69, 7, 7),
    @(70, 8, 8),
    @(71, 7, 8),
    @(72, 8, 8),
    @(73, 2, 2),
    @(74, 7, 7),
    @(75, 9, 9),
    @(76, 9, 9),
    @(77, 9, 9),
    @(78, 4, 4),
    @(79, 8, 8),
    @(80, 6, 6),
    @(81, 4, 4),
    @(82, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iv = $row[1]
    $jv = $row[2]
    $ws.Cells.Item($r, 9).Value = $iv
    $ws.Cells.Item($r, 10).Value = $jv
}
